$d = $word.ActiveDocument

# 1) Remove the trailing space after the Jersey Dependency sentence.
$d.Content.Find.Execute(
    "JavaFX client application. ", $false, $false, $false, $false, $false,
    $true, 1, $false, "JavaFX client application.", 2) | Out-Null

# 2) Find the paragraph that now ends with the Jersey Dependency sentence and
#    insert a brand-new list paragraph right after it (InsertParagraphAfter
#    clones the paragraph's pPr/rPr, including the blue-underline list
#    formatting, into a fresh empty paragraph).
$targetIdx = -1
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "*Jersey Dependency*JavaFX client application.*") {
        $targetIdx = $idx
    }
}

$p = $d.Paragraphs.Item($targetIdx)
$p.Range.InsertParagraphAfter() | Out-Null

$newPara = $d.Paragraphs.Item($targetIdx + 1)
$newRange = $newPara.Range
$insertStart = $newRange.Start

# 3) Fill the new paragraph with the Apache Tomcat sentence (plus trailing
#    space), then drop a collapsed "_GoBack" bookmark between "w" and "eb".
$newRange.InsertBefore("Apache Tomcat version 8.0.52 as a web application. ")

$splitPos = $insertStart + "Apache Tomcat version 8.0.52 as a w".Length
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# 4) Merge the old "After system " / "work" split run back into a single
#    run; replacing the visible text across the bookmarked boundary also
#    removes the stale "_GoBack" bookmark that used to sit there.
$d.Content.Find.Execute(
    "After system work", $false, $false, $false, $false, $false,
    $true, 1, $false, "After system work", 2) | Out-Null
